$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.826.80'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '3.085.30'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').Value = '595.62'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').Value = '155.96'
$ws.Range('E6').Value = '  +1.91%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = '0.538'
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').Value = '3.080.82'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('E10').Value = '  -2.42%  '
$ws.Range('D11').Value = '5.95'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('D13').Value = '0.0000238'
$ws.Range('E13').Value = '  -3.82%  '
$ws.Range('D14').Value = '36.80'
$ws.Range('E14').Value = '  -4.47%  '
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '3.597.90'
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = '7.20'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').Value = '63.787.14'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('D19').Value = '3.092.71'
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('D20').Value = '479.28'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = '14.45'
$ws.Range('E21').Value = '  -3.55%  '
$ws.Range('D22').Value = '0.712'
$ws.Range('E22').Value = '  -4.69%  '
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('D24').Value = '2.47'
$ws.Range('E24').Value = '  +2.59%  '
$ws.Range('D25').Value = '81.53'
$ws.Range('D26').Value = '12.86'
$ws.Range('E26').Value = '  -4.86%  '
$ws.Range('D27').Value = '10.66'
$ws.Range('E27').Value = '  +8.16%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').Value = '7.64'
$ws.Range('E29').Value = '  +2.44%  '
$ws.Range('D30').Value = '2.69'
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').Value = '2.20'
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('D33').Value = '0.112'
$ws.Range('E33').Value = '  -4.80%  '
$ws.Range('D34').Value = '27.16'
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('D35').Value = '0.0₃0840'
$ws.Range('E35').Value = '  -4.17%  '
$ws.Range('D36').Value = '1.07'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').Value = '6.02'
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  -6.20%  '
$ws.Range('D39').Value = '2.24'
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('D40').Value = '50.90'
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').Value = '441.67'
$ws.Range('E42').Value = '  -5.53%  '
$ws.Range('D43').Value = '0.289'
$ws.Range('E43').Value = '  -3.94%  '
$ws.Range('D44').Value = '0.0364'
$ws.Range('E44').Value = '  -4.42%  '
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').Value = '2.834.91'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').Value = '40.05'
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('D48').Value = '132.24'
$ws.Range('E48').Value = '  +0.58%  '
$ws.Range('D49').Value = '25.95'
$ws.Range('E49').Value = '  +0.65%  '
$ws.Range('D51').Value = '2.24'
$ws.Range('E51').Value = '  -2.15%  '
